$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source row for case 4768 ("VALLESE, FELIPE 684", old row 24) was
# removed from the tracking sheet. Deleting the entire row shifts every
# subsequent record up by one position and drops the table's final row
# (old row 90, case 6917) from the used range, matching the new
# dimension of A1:P89.
$ws.Rows.Item(24).Delete()
